$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.422.06"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "2.218.88"
$ws.Range("E3").Value = "  -1.30%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "297.96"
$ws.Range("E5").Value = "  -3.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "90.30"
$ws.Range("E6").Value = "  -4.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.565"
$ws.Range("E7").Value = "  -1.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.01"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.491"
$ws.Range("E9").Value = "  -6.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.14"
$ws.Range("E10").Value = "  -6.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0787"
$ws.Range("E11").Value = "  -3.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.97"
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("D14").Value = "2.557.54"
$ws.Range("E14").Value = "  -1.30%  "
$ws.Range("D15").Value = "2.264.46"
$ws.Range("E15").Value = "  -6.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.53"
$ws.Range("E16").Value = "  -1.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.788"
$ws.Range("E17").Value = "  -6.29%  "
$ws.Range("D18").Value = "44.169.87"
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("D19").Value = "0.0₃0906"
$ws.Range("E19").Value = "  -6.38%  "
$ws.Range("E20").Value = "  -7.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.12"
$ws.Range("E21").Value = "  -9.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "64.29"
$ws.Range("E22").Value = "  -2.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.27"
$ws.Range("E23").Value = "  -0.66%  "
$ws.Range("E24").Value = "  -12.84%  "
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("E26").Value = "  -6.97%  "
$ws.Range("E27").Value = "  +0.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.52"
$ws.Range("E28").Value = "  +0.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.38"
$ws.Range("E29").Value = "  -4.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.54"
$ws.Range("E30").Value = "  -2.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "148.61"
$ws.Range("E31").Value = "  -3.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.39"
$ws.Range("E32").Value = "  -10.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.54"
$ws.Range("E33").Value = "  -3.40%  "
$ws.Range("E34").Value = "  -6.08%  "
$ws.Range("E35").Value = "  -3.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.86"
$ws.Range("E36").Value = "  -8.16%  "
$ws.Range("E37").Value = "  -4.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.67"
$ws.Range("E38").Value = "  -7.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0298"
$ws.Range("E39").Value = "  -2.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.20"
$ws.Range("E40").Value = "  -7.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.57"
$ws.Range("E41").Value = "  -7.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.23"
$ws.Range("E42").Value = "  -9.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.01"
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("D44").Value = "1.824.21"
$ws.Range("E44").Value = "  +4.12%  "
$ws.Range("E45").Value = "  +10.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.179"
$ws.Range("E46").Value = "  -7.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "67.57"
$ws.Range("E47").Value = "  -5.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "73.94"
$ws.Range("E48").Value = "  -8.68%  "
$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.84"
$ws.Range("E49").Value = "  +14.18%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "94.47"
$ws.Range("E50").Value = "  -5.53%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.439.73"
$ws.Range("E51").Value = "  -1.24%  "
